$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for Alcachofa / Macroferia Regional de
# Talca. It belongs right after the current header block of recent rows, so
# insert a fresh row at 21 and push the existing data (rows 21:73) down to
# 22:74 — exactly like typing a new row into the middle of the table in
# Excel.
$ws.Rows("21:21").Insert()

# Populate the newly inserted row with the new record's values.
$ws.Range("A21").Value = 5
$ws.Range("B21").Value = "Macroferia Regional de Talca"
$ws.Range("C21").Value = "Maule"
$ws.Range("D21").Value = 44498
$ws.Range("E21").Value = 7
$ws.Range("F21").Value = 100112013
$ws.Range("G21").Value = "Alcachofa"
$ws.Range("H21").Value = "Madrigal"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = 10000
$ws.Range("N21").Value = "`$/caja 40 unidades"
$ws.Range("O21").Value = "Provincia del Elquí"
$ws.Range("P21").Value = 250
$ws.Range("Q21").Value = 40
$ws.Range("R21").Value = "Hortaliza"
